# Adapt column header formatting to respective input file names (FV2404 / FV2410)
# and turn the header row + data range into a proper Excel Table, with the
# header row frozen in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" / "_new" suffixed headers to "_FV2404" / "_FV2410".
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($c in $oldCols) {
    $cell = $ws.Range($c + "1")
    $v = $cell.Value()
    $cell.Value = $v -replace "_old$", "_FV2404"
}
foreach ($c in $newCols) {
    $cell = $ws.Range($c + "1")
    $v = $cell.Value()
    $cell.Value = $v -replace "_new$", "_FV2410"
}

# 2) Turn the used range into an Excel Table ("Table1") with an AutoFilter.
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
